$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 87.52000886463766
$ws.Range("C2").Value = 122.3238612539563
$ws.Range("D2").Value = 141.4310488244865
$ws.Range("E2").Value = 154.0928370029275

$ws.Range("B3").Value = 108.6750916094591
$ws.Range("C3").Value = 151.0257897479978
$ws.Range("D3").Value = 171.0166763028857
$ws.Range("E3").Value = 185.9038688915234

$ws.Range("B4").Value = 88.89773292123303
$ws.Range("C4").Value = 126.4672401468646
$ws.Range("D4").Value = 147.9905732170969
$ws.Range("E4").Value = 165.5804725489884

$ws.Range("B5").Value = 77.18715505475731
$ws.Range("C5").Value = 106.3103909170766
$ws.Range("D5").Value = 116.8595139535432
$ws.Range("E5").Value = 125.956289484215

$ws.Range("B6").Value = 67.25326813445629
$ws.Range("C6").Value = 92.26178911390346
$ws.Range("D6").Value = 102.3139385843235
$ws.Range("E6").Value = 109.2804573964136

$ws.Range("B7").Value = 7.400791066150671
$ws.Range("C7").Value = 10.071140552069
$ws.Range("D7").Value = 11.1905147623712
$ws.Range("E7").Value = 11.79948541821851

$ws.Range("B8").Value = 353.1640303372616
$ws.Range("C8").Value = 493.2740685376367
$ws.Range("D8").Value = 565.1268577930598
$ws.Range("E8").Value = 602.5235190325914

$ws.Range("B9").Value = 102.4334863009774
$ws.Range("C9").Value = 140.362942969836
$ws.Range("D9").Value = 155.143863274867
$ws.Range("E9").Value = 164.274537571003

$ws.Range("B10").Value = 44.46686990261082
$ws.Range("C10").Value = 59.02903619500883
$ws.Range("D10").Value = 65.55865434390834
$ws.Range("E10").Value = 67.50467295181468

$ws.Range("B11").Value = 8.101745782040659
$ws.Range("C11").Value = 10.19151698997291
$ws.Range("D11").Value = 11.18548547597686
$ws.Range("E11").Value = 12.43537387220275

$ws.Range("B12").Value = 18.44521984358384
$ws.Range("C12").Value = 24.62685994101061
$ws.Range("D12").Value = 26.16350002395977
$ws.Range("E12").Value = 26.13346335982067

$ws.Range("B13").Value = 24.93231680160065
$ws.Range("C13").Value = 32.76012645591509
$ws.Range("D13").Value = 36.65814988539523
$ws.Range("E13").Value = 38.16989894895444

